$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column cells are stored as text, preserving exact formatting
# (trailing zeros, thousand-dot separators, etc.) as in the source data,
# instead of Excel auto-converting numeric-looking strings to numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "45.019.28"
$ws.Range("E2").Value = "  +0.76%  "
$ws.Range("D3").Value = "2.265.22"
$ws.Range("E3").Value = "  +0.89%  "
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  -0.65%  "
$ws.Range("D5").Value = "301.65"
$ws.Range("E5").Value = "  -1.70%  "
$ws.Range("D6").Value = "94.40"
$ws.Range("E6").Value = "  -1.91%  "
$ws.Range("D7").Value = "0.565"
$ws.Range("E7").Value = "  -1.21%  "
$ws.Range("E8").Value = "  -0.52%  "
$ws.Range("D9").Value = "0.509"
$ws.Range("E9").Value = "  -2.72%  "
$ws.Range("E10").Value = "  -2.56%  "
$ws.Range("E11").Value = "  -2.26%  "
$ws.Range("D12").Value = "7.21"
$ws.Range("E12").Value = "  -0.89%  "
$ws.Range("E13").Value = "  -0.56%  "
$ws.Range("D14").Value = "2.612.26"
$ws.Range("E14").Value = "  +1.00%  "
$ws.Range("D15").Value = "2.268.35"
$ws.Range("E15").Value = "  -1.56%  "
$ws.Range("D16").Value = "13.76"
$ws.Range("E16").Value = "  +0.75%  "
$ws.Range("D17").Value = "0.796"
$ws.Range("E17").Value = "  -5.26%  "
$ws.Range("D18").Value = "44.867.60"
$ws.Range("E18").Value = "  +1.05%  "
$ws.Range("D19").Value = "12.91"
$ws.Range("E19").Value = "  +7.34%  "
$ws.Range("E20").Value = "  -3.34%  "
$ws.Range("D21").Value = "6.08"
$ws.Range("E21").Value = "  -3.81%  "
$ws.Range("D22").Value = "65.23"
$ws.Range("E22").Value = "  -0.58%  "
$ws.Range("D23").Value = "238.79"
$ws.Range("E23").Value = "  +0.09%  "
$ws.Range("E24").Value = "  -2.77%  "
$ws.Range("E25").Value = "  -0.35%  "
$ws.Range("D26").Value = "1.92"
$ws.Range("E26").Value = "  -3.77%  "
$ws.Range("D27").Value = "41.95"
$ws.Range("E27").Value = "  +12.17%  "
$ws.Range("E28").Value = "  +0.81%  "
$ws.Range("D29").Value = "9.54"
$ws.Range("E29").Value = "  -3.14%  "
$ws.Range("D30").Value = "19.53"
$ws.Range("E30").Value = "  -2.40%  "
$ws.Range("D31").Value = "152.54"
$ws.Range("E31").Value = "  -0.07%  "
$ws.Range("D32").Value = "5.61"
$ws.Range("E32").Value = "  -6.75%  "
$ws.Range("D33").Value = "0.0785"
$ws.Range("E33").Value = "  -1.99%  "
$ws.Range("D34").Value = "2.56"
$ws.Range("E34").Value = "  -3.25%  "
$ws.Range("B35").Value = "LidoDAOToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D35").Value = "2.92"
$ws.Range("E35").Value = "  -4.30%  "
$ws.Range("B36").Value = "Stellar"
$ws.Range("C36").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D36").Value = "0.117"
$ws.Range("E36").Value = "  -1.68%  "
$ws.Range("D37").Value = "0.105"
$ws.Range("E37").Value = "  -3.69%  "
$ws.Range("D38").Value = "1.73"
$ws.Range("E38").Value = "  -6.72%  "
$ws.Range("E39").Value = "  +0.98%  "
$ws.Range("B40").Value = "RenderToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D40").Value = "3.80"
$ws.Range("E40").Value = "  +1.23%  "
$ws.Range("B41").Value = "NEARProtocol"
$ws.Range("C41").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D41").Value = "3.26"
$ws.Range("E41").Value = "  -4.41%  "
$ws.Range("D42").Value = "13.82"
$ws.Range("E42").Value = "  -9.67%  "
$ws.Range("D43").Value = "1.00"
$ws.Range("E43").Value = "  -0.61%  "
$ws.Range("B44").Value = "Stacks"
$ws.Range("C44").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D44").Value = "1.94"
$ws.Range("E44").Value = "  +13.10%  "
$ws.Range("B45").Value = "Maker"
$ws.Range("C45").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D45").Value = "1.782.29"
$ws.Range("E45").Value = "  -2.64%  "
$ws.Range("D46").Value = "0.192"
$ws.Range("E46").Value = "  -0.85%  "
$ws.Range("D47").Value = "69.96"
$ws.Range("E47").Value = "  -0.89%  "
$ws.Range("D48").Value = "75.46"
$ws.Range("E48").Value = "  -5.23%  "
$ws.Range("D49").Value = "96.77"
$ws.Range("E49").Value = "  -2.69%  "
$ws.Range("B50").Value = "FraxShare"
$ws.Range("C50").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D50").Value = "7.86"
$ws.Range("E50").Value = "  -2.53%  "
$ws.Range("B51").Value = "MultiversX"
$ws.Range("C51").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D51").Value = "52.98"
$ws.Range("E51").Value = "  -3.23%  "
